$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Timetable data grid (rows 4-12, columns B-U) with the
# reshuffled class-slot assignments.
$ws.Range("C4").Value = "B Tech IT"
$ws.Range("D4").Value = "MBA Tech CS"
$ws.Range("E4").Value = "MBA Tech IT"
$ws.Range("H4").Value = "MBA Tech IT"
$ws.Range("I4").Value = "MBA Tech CS"
$ws.Range("J4").Value = "MBA Tech CS"
$ws.Range("K4").Value = "MBA Tech IT"
$ws.Range("L4").Value = "B Tech CS"
$ws.Range("M4").Value = "B Tech IT"
$ws.Range("N4").Value = "MBA Tech IT"
$ws.Range("O4").Value = "MBA Tech CS"
$ws.Range("P4").Value = "B Tech IT"
$ws.Range("Q4").Value = "B Tech CS"
$ws.Range("R4").Value = "B Tech CS"
$ws.Range("S4").Value = "B Tech IT"
$ws.Range("T4").Value = "MBA Tech CS"
$ws.Range("U4").Value = "MBA Tech IT"
$ws.Range("B5").Value = "B Tech IT"
$ws.Range("D5").Value = "MBA Tech IT"
$ws.Range("E5").Value = "MBA Tech CS"
$ws.Range("F5").Value = "MBA Tech CS"
$ws.Range("H5").Value = "B Tech CS"
$ws.Range("I5").Value = "B Tech IT"
$ws.Range("J5").Value = "MBA Tech IT"
$ws.Range("K5").Value = "MBA Tech CS"
$ws.Range("L5").Value = "B Tech IT"
$ws.Range("N5").Value = "B Tech CS"
$ws.Range("O5").Value = "B Tech IT"
$ws.Range("P5").Value = "MBA Tech CS"
$ws.Range("Q5").Value = "MBA Tech IT"
$ws.Range("R5").Value = "B Tech IT"
$ws.Range("S5").Value = "B Tech CS"
$ws.Range("T5").Value = "MBA Tech IT"
$ws.Range("U5").Value = "MBA Tech CS"
$ws.Range("B6").Value = "MBA Tech CS"
$ws.Range("C6").Value = "MBA Tech IT"
$ws.Range("D6").Value = "B Tech CS"
$ws.Range("E6").Value = "B Tech IT"
$ws.Range("G6").Value = "MBA Tech CS"
$ws.Range("H6").Value = "B Tech IT"
$ws.Range("I6").Value = "B Tech CS"
$ws.Range("K6").Value = "B Tech IT"
$ws.Range("L6").Value = "MBA Tech CS"
$ws.Range("M6").Value = "MBA Tech IT"
$ws.Range("O6").Value = "B Tech CS"
$ws.Range("P6").Value = "MBA Tech IT"
$ws.Range("Q6").Value = "MBA Tech CS"
$ws.Range("R6").Value = "MBA Tech CS"
$ws.Range("S6").Value = "MBA Tech IT"
$ws.Range("T6").Value = "B Tech CS"
$ws.Range("U6").Value = "B Tech IT"
$ws.Range("B7").Value = "MBA Tech IT"
$ws.Range("C7").Value = "MBA Tech CS"
$ws.Range("D7").Value = "B Tech IT"
$ws.Range("E7").Value = "B Tech CS"
$ws.Range("F7").Value = "B Tech CS"
$ws.Range("G7").Value = "B Tech IT"
$ws.Range("H7").Value = "MBA Tech CS"
$ws.Range("I7").Value = "MBA Tech IT"
$ws.Range("J7").Value = "B Tech IT"
$ws.Range("K7").Value = "B Tech CS"
$ws.Range("L7").Value = "MBA Tech IT"
$ws.Range("M7").Value = "MBA Tech CS"
$ws.Range("N7").Value = "MBA Tech CS"
$ws.Range("O7").Value = "MBA Tech IT"
$ws.Range("P7").Value = "B Tech CS"
$ws.Range("Q7").Value = "B Tech IT"
$ws.Range("R7").Value = "MBA Tech IT"
$ws.Range("S7").Value = "MBA Tech CS"
$ws.Range("T7").Value = "B Tech IT"
$ws.Range("U7").Value = "B Tech CS"
$ws.Range("B9").Value = "B Tech CS"
$ws.Range("C9").Value = "B Tech IT"
$ws.Range("D9").Value = "MBA Tech CS"
$ws.Range("E9").Value = "MBA Tech IT"
$ws.Range("F9").Value = "B Tech IT"
$ws.Range("G9").Value = "B Tech CS"
$ws.Range("I9").Value = "MBA Tech CS"
$ws.Range("J9").Value = "MBA Tech CS"
$ws.Range("K9").Value = "MBA Tech IT"
$ws.Range("L9").Value = "B Tech CS"
$ws.Range("M9").Value = "B Tech IT"
$ws.Range("N9").Value = "MBA Tech IT"
$ws.Range("O9").Value = "MBA Tech CS"
$ws.Range("P9").Value = "B Tech IT"
$ws.Range("Q9").Value = "B Tech CS"
$ws.Range("R9").Value = "B Tech CS"
$ws.Range("S9").Value = "B Tech IT"
$ws.Range("U9").Value = "MBA Tech IT"
$ws.Range("B10").Value = "B Tech IT"
$ws.Range("C10").Value = "B Tech CS"
$ws.Range("D10").Value = "MBA Tech IT"
$ws.Range("E10").Value = "MBA Tech CS"
$ws.Range("F10").Value = "MBA Tech CS"
$ws.Range("G10").Value = "MBA Tech IT"
$ws.Range("H10").Value = "B Tech CS"
$ws.Range("I10").Value = "B Tech IT"
$ws.Range("K10").Value = "MBA Tech CS"
$ws.Range("L10").Value = "B Tech IT"
$ws.Range("M10").Value = "B Tech CS"
$ws.Range("N10").Value = "B Tech CS"
$ws.Range("O10").Value = "B Tech IT"
$ws.Range("P10").Value = "MBA Tech CS"
$ws.Range("Q10").Value = "MBA Tech IT"
$ws.Range("R10").Value = "B Tech IT"
$ws.Range("S10").Value = "B Tech CS"
$ws.Range("T10").Value = "MBA Tech IT"
$ws.Range("B11").Value = "MBA Tech CS"
$ws.Range("C11").Value = "MBA Tech IT"
$ws.Range("D11").Value = "B Tech CS"
$ws.Range("E11").Value = "B Tech IT"
$ws.Range("G11").Value = "MBA Tech CS"
$ws.Range("H11").Value = "B Tech IT"
$ws.Range("I11").Value = "B Tech CS"
$ws.Range("J11").Value = "B Tech CS"
$ws.Range("K11").Value = "B Tech IT"
$ws.Range("L11").Value = "MBA Tech CS"
$ws.Range("M11").Value = "MBA Tech IT"
$ws.Range("N11").Value = "B Tech IT"
$ws.Range("O11").Value = "B Tech CS"
$ws.Range("P11").Value = "MBA Tech IT"
$ws.Range("Q11").Value = "MBA Tech CS"
$ws.Range("R11").Value = "MBA Tech CS"
$ws.Range("S11").Value = "MBA Tech IT"
$ws.Range("T11").Value = "B Tech CS"
$ws.Range("U11").Value = "B Tech IT"
$ws.Range("B12").Value = "MBA Tech IT"
$ws.Range("C12").Value = "MBA Tech CS"
$ws.Range("D12").Value = "B Tech IT"
$ws.Range("E12").Value = "B Tech CS"
$ws.Range("F12").Value = "B Tech CS"
$ws.Range("H12").Value = ""
$ws.Range("J12").Value = "B Tech CS"
$ws.Range("L12").Value = ""
$ws.Range("N12").Value = "B Tech CS"
$ws.Range("P12").Value = ""
$ws.Range("R12").Value = "B Tech CS"
$ws.Range("T12").Value = ""
$ws.Range("U12").Value = ""
